# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on rows 3 and 5
# (the db324f70... handoff/handback entry) for both the zh-cn and
# de-de worksheets, reflecting a newly generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-24 02:23:09"
$wsZhCn.Range("H3").Value = "2016-03-24 02:23:35"
$wsZhCn.Range("E5").Value = "2016-03-24 02:23:09"
$wsZhCn.Range("H5").Value = "2016-03-24 02:23:35"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-24 02:23:13"
$wsDeDe.Range("H3").Value = "2016-03-24 02:23:42"
$wsDeDe.Range("E5").Value = "2016-03-24 02:23:13"
$wsDeDe.Range("H5").Value = "2016-03-24 02:23:42"
